$wb = $excel.ActiveWorkbook

# 1. Update the "Date" value on the Metadata sheet
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# 2. Update the System URI for CodeSystem-TRE-R288-TypeProfession on "Include #0"
$inc0Sheet = $wb.Worksheets.Item("Include #0")
$inc0Sheet.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R288-TypeProfession"

# 3. Update the System URI for CodeSystem-TRE-R37-TypeProfessionFonction on "Include #1"
$inc1Sheet = $wb.Worksheets.Item("Include #1")
$inc1Sheet.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R37-TypeProfessionFonction"

$wb.Save()
